$d = $word.ActiveDocument

$old = "Supported quality improvements by contributing to design reviews, communicating test results to stakeholders, and maintaining automation scripts and frameworks."
$new = "Supported quality improvements by contributing to design reviews & communicating test results to stakeholders."

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
